# Root Chord Results.xlsx - "Parametric Study" sheet update
#
# The parametric sweep table (rows 3-5, columns B:CT) had its three input
# rows rotated by one position:
#   new row 3 = old row 5   ("Root Chord" row moves to the top)
#   new row 4 = old row 3   ("Mach number" row moves down one)
#   new row 5 = old row 4   ("Angle of attack" row moves down one)
# and the sheet's active cell selection moved from B11 to E35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 2   # column B
$lastCol  = 98  # column CT

# Snapshot the current contents of rows 3, 4 and 5 (labels in col B are
# shared-string text, the rest are numeric - Value2 round-trips both).
$row3 = @{}
$row4 = @{}
$row5 = @{}

for ($col = $firstCol; $col -le $lastCol; $col++) {
    $row3[$col] = $ws.Cells.Item(3, $col).Value2
    $row4[$col] = $ws.Cells.Item(4, $col).Value2
    $row5[$col] = $ws.Cells.Item(5, $col).Value2
}

# Write the rotated values back: row3<-row5, row4<-row3(old), row5<-row4(old)
for ($col = $firstCol; $col -le $lastCol; $col++) {
    $ws.Cells.Item(3, $col).Value = $row5[$col]
    $ws.Cells.Item(4, $col).Value = $row3[$col]
    $ws.Cells.Item(5, $col).Value = $row4[$col]
}

# Move the selection to match the saved view state (E35).
[void]$ws.Range("E35").Select()
